$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (email column): "1234sals!" -> "1234sals!s"
$ws.Range("B7").Value = "1234sals!s"

# Row 8 (email column): "sals12@gmail.com" -> "sals1234@gmail.com"
$ws.Range("B8").Value = "sals1234@gmail.com"

# Move the active selection from G7 to G10
$ws.Range("G10").Select() | Out-Null
